$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1017.3333
$ws.Cells.Item(28, 9).Value = 753.3684
$ws.Cells.Item(28, 10).Value = 2020.4
$ws.Cells.Item(28, 11).Value = 753.3684
$ws.Cells.Item(28, 12).Value = 2020.4
$ws.Cells.Item(28, 13).Value = -268.3684
$ws.Cells.Item(28, 14).Value = -2990.4
$ws.Cells.Item(40, 8).Value = 1797.1666
$ws.Cells.Item(40, 9).Value = 1556.6
$ws.Cells.Item(40, 10).Value = 3000
$ws.Cells.Item(40, 11).Value = 1556.6
$ws.Cells.Item(40, 12).Value = 3000
$ws.Cells.Item(40, 13).Value = -1381.6
$ws.Cells.Item(40, 14).Value = -3350
$ws.Cells.Item(55, 8).Value = 145.45454
$ws.Cells.Item(55, 9).Value = 120
$ws.Cells.Item(55, 11).Value = 120
$ws.Cells.Item(55, 13).Value = 94
$ws.Cells.Item(62, 8).Value = 2500
$ws.Cells.Item(62, 9).Value = 2500
$ws.Cells.Item(62, 10).Value = 2500
$ws.Cells.Item(62, 11).Value = 2500
$ws.Cells.Item(62, 12).Value = 2500
$ws.Cells.Item(62, 13).Value = -1876
$ws.Cells.Item(62, 14).Value = -3748
$ws.Cells.Item(65, 8).Value = 2500
$ws.Cells.Item(65, 9).Value = 2500
$ws.Cells.Item(65, 10).Value = 2500
$ws.Cells.Item(65, 11).Value = 12500
$ws.Cells.Item(65, 12).Value = 12500
$ws.Cells.Item(65, 13).Value = -9380
$ws.Cells.Item(65, 14).Value = -18740
$ws.Cells.Item(116, 8).Value = 7201.4443
$ws.Cells.Item(116, 10).Value = 3264
$ws.Cells.Item(116, 12).Value = 3264
$ws.Cells.Item(116, 14).Value = -10148
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1112.82
$ws.Cells.Item(32, 9).Value = 860.525
$ws.Cells.Item(32, 11).Value = 860.525
$ws.Cells.Item(32, 13).Value = -573.525
$ws.Cells.Item(34, 8).Value = 6000
$ws.Cells.Item(34, 10).Value = 6000
$ws.Cells.Item(34, 12).Value = 6000
$ws.Cells.Item(34, 14).Value = -6542
$ws.Cells.Item(45, 8).Value = 3301.7693
$ws.Cells.Item(45, 9).Value = 3155.6365
$ws.Cells.Item(45, 10).Value = 4105.5
$ws.Cells.Item(45, 11).Value = 3155.6365
$ws.Cells.Item(45, 12).Value = 4105.5
$ws.Cells.Item(45, 13).Value = -2778.6365
$ws.Cells.Item(45, 14).Value = -4859.5
$ws.Cells.Item(122, 8).Value = 1956.7368
$ws.Cells.Item(122, 9).Value = 1684.909
$ws.Cells.Item(122, 10).Value = 2330.5
$ws.Cells.Item(122, 11).Value = 5054.727000000001
$ws.Cells.Item(122, 12).Value = 6991.5
$ws.Cells.Item(122, 13).Value = -2604.727000000001
$ws.Cells.Item(122, 14).Value = -11891.5
$ws.Cells.Item(132, 8).Value = 6237
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 6237
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).ClearContents()
$ws.Cells.Item(132, 13).Value = 18711
$ws.Cells.Item(132, 14).Value = -23771
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 3640.889
$ws.Cells.Item(36, 9).Value = 3640.889
$ws.Cells.Item(36, 11).Value = 3640.889
$ws.Cells.Item(36, 13).Value = -3106.889
$ws.Cells.Item(86, 8).Value = 10865.467
$ws.Cells.Item(86, 9).Value = 4175.5835
$ws.Cells.Item(86, 10).Value = 37625
$ws.Cells.Item(86, 11).Value = 4175.5835
$ws.Cells.Item(86, 12).Value = 37625
$ws.Cells.Item(86, 13).Value = -3052.5835
$ws.Cells.Item(86, 14).Value = -39871
$ws.Cells.Item(89, 8).Value = 10865.467
$ws.Cells.Item(89, 9).Value = 4175.5835
$ws.Cells.Item(89, 10).Value = 37625
$ws.Cells.Item(89, 11).Value = 20877.9175
$ws.Cells.Item(89, 12).Value = 188125
$ws.Cells.Item(89, 13).Value = -15261.9175
$ws.Cells.Item(89, 14).Value = -199357
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(60, 8).Value = 11442.5
$ws.Cells.Item(74, 8).Value = 39990
$ws.Cells.Item(74, 10).Value = 39990
$ws.Cells.Item(74, 12).Value = 39990
$ws.Cells.Item(74, 14).Value = -41738
$ws.Cells.Item(77, 8).Value = 39990
$ws.Cells.Item(77, 10).Value = 39990
$ws.Cells.Item(77, 12).Value = 119970
$ws.Cells.Item(77, 14).Value = -128706
$ws.Cells.Item(107, 8).Value = 2742.8696
$ws.Cells.Item(107, 9).Value = 2917.4
$ws.Cells.Item(107, 11).Value = 2917.4
$ws.Cells.Item(107, 13).Value = -997.4000000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 335.2353
$ws.Cells.Item(12, 9).Value = 9.125
$ws.Cells.Item(12, 10).Value = 625.1111
$ws.Cells.Item(12, 11).Value = 27.375
$ws.Cells.Item(12, 12).Value = 1875.3333
$ws.Cells.Item(12, 13).Value = 145.625
$ws.Cells.Item(12, 14).Value = -2221.3333
$ws.Cells.Item(51, 8).Value = 1999.5
$ws.Cells.Item(51, 10).Value = 1999
$ws.Cells.Item(51, 12).Value = 5997
$ws.Cells.Item(51, 14).Value = -6917
$ws.Cells.Item(100, 8).Value = 5413.3335
$ws.Cells.Item(107, 8).Value = 1211.317
$ws.Cells.Item(107, 9).Value = 255
$ws.Cells.Item(107, 10).Value = 1375.2572
$ws.Cells.Item(107, 11).Value = 765
$ws.Cells.Item(107, 12).Value = 4125.7716
$ws.Cells.Item(107, 13).Value = 1155
$ws.Cells.Item(107, 14).Value = -7965.7716
$ws.Cells.Item(122, 8).Value = 467
$ws.Cells.Item(122, 9).Value = 881.75
$ws.Cells.Item(122, 10).Value = 316.18182
$ws.Cells.Item(122, 11).Value = 7935.75
$ws.Cells.Item(122, 12).Value = 2845.63638
$ws.Cells.Item(122, 13).Value = -5485.75
$ws.Cells.Item(122, 14).Value = -7745.63638
$ws.Cells.Item(137, 8).Value = 30000
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 30000
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).ClearContents()
$ws.Cells.Item(137, 13).Value = 90000
$ws.Cells.Item(137, 14).Value = -100200
$ws.Cells.Item(139, 8).Value = 20629.555
$ws.Cells.Item(139, 9).Value = 2000
$ws.Cells.Item(139, 10).Value = 22958.25
$ws.Cells.Item(139, 11).Value = 6000
$ws.Cells.Item(139, 12).Value = 68874.75
$ws.Cells.Item(139, 13).Value = -860
$ws.Cells.Item(139, 14).Value = -79154.75
$ws.Cells.Item(141, 8).Value = 42997.35
$ws.Cells.Item(141, 10).Value = 50000
$ws.Cells.Item(141, 12).Value = 150000
$ws.Cells.Item(141, 14).Value = -160360
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2076.077
$ws.Cells.Item(113, 9).Value = 1599.1
$ws.Cells.Item(113, 11).Value = 1599.1
$ws.Cells.Item(113, 13).Value = 570.9000000000001
$ws.Cells.Item(132, 8).Value = 2397.2415
$ws.Cells.Item(132, 9).Value = 2454.4583
$ws.Cells.Item(132, 10).Value = 2122.6
$ws.Cells.Item(132, 11).Value = 7363.374899999999
$ws.Cells.Item(132, 12).Value = 6367.799999999999
$ws.Cells.Item(132, 13).Value = -4833.374899999999
$ws.Cells.Item(132, 14).Value = -11427.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2133.1667
$ws.Cells.Item(46, 10).Value = 3166.3333
$ws.Cells.Item(46, 12).Value = 3166.3333
$ws.Cells.Item(46, 14).Value = -3542.3333
$ws.Cells.Item(122, 8).Value = 4207.3687
$ws.Cells.Item(122, 9).Value = 4294
$ws.Cells.Item(122, 11).Value = 12882
$ws.Cells.Item(122, 13).Value = -10432
$ws.Cells.Item(130, 8).Value = 89999
$ws.Cells.Item(130, 10).Value = 89999
$ws.Cells.Item(130, 12).Value = 89999
$ws.Cells.Item(130, 14).Value = -100039
$ws.Cells.Item(132, 8).Value = 2227.8386
$ws.Cells.Item(132, 9).Value = 1218.579
$ws.Cells.Item(132, 11).Value = 3655.737
$ws.Cells.Item(132, 13).Value = -1125.737
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 53962.555
$ws.Cells.Item(41, 9).Value = 50331
$ws.Cells.Item(41, 10).Value = 55778.332
$ws.Cells.Item(41, 11).Value = 50331
$ws.Cells.Item(41, 12).Value = 55778.332
$ws.Cells.Item(41, 13).Value = -49941
$ws.Cells.Item(41, 14).Value = -56558.332
$ws.Cells.Item(132, 8).Value = 17546410
$ws.Cells.Item(132, 9).Value = 18184096
$ws.Cells.Item(132, 11).Value = 54552288
$ws.Cells.Item(132, 13).Value = -54549758
